# Applies the "Updated cryptos list" refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $range while forcing a text (string) result,
# even when $text looks like a number (e.g. "605.82"). A direct
# `.Value = "605.82"` assignment gets auto-coerced to a numeric cell by
# the COM layer's type inference, but the source workbook stores these as
# plain inline strings (t="inlineStr"), so we round-trip through a helper
# cell holding a text-returning formula and paste-special the computed
# *value* (xlPasteValues = -4163), which preserves the text type exactly
# like a manual Paste-Values would in real Excel.
function Set-TextValue($range, [string]$text) {
    $helperCell = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $helperCell.Formula = '="' + $escaped + '"'
    $helperCell.Copy()
    $range.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "63.377.05"
$ws.Range("E2").Value = "  -2.60%  "

$ws.Range("D3").Value = "3.112.17"
$ws.Range("E3").Value = "  -2.98%  "

$ws.Range("E4").Value = "  +0.22%  "

Set-TextValue $ws.Range("D5") "605.82"
$ws.Range("E5").Value = "  +0.61%  "

Set-TextValue $ws.Range("D6") "145.17"
$ws.Range("E6").Value = "  -5.75%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "3.105.05"
$ws.Range("E8").Value = "  -3.23%  "

$ws.Range("E9").Value = "  -2.21%  "

Set-TextValue $ws.Range("D10") "0.149"
$ws.Range("E10").Value = "  -3.83%  "

Set-TextValue $ws.Range("D11") "5.27"
$ws.Range("E11").Value = "  -5.26%  "

Set-TextValue $ws.Range("D12") "0.467"
$ws.Range("E12").Value = "  -3.46%  "

Set-TextValue $ws.Range("D13") "0.0000250"
$ws.Range("E13").Value = "  -3.45%  "

Set-TextValue $ws.Range("D14") "35.09"
$ws.Range("E14").Value = "  -6.43%  "

$ws.Range("D15").Value = "3.626.57"
$ws.Range("E15").Value = "  -2.86%  "

$ws.Range("E16").Value = "  +2.23%  "

$ws.Range("D17").Value = "63.440.26"
$ws.Range("E17").Value = "  -2.71%  "

$ws.Range("D18").Value = "3.116.35"
$ws.Range("E18").Value = "  -2.89%  "

Set-TextValue $ws.Range("D19") "6.79"
$ws.Range("E19").Value = "  -4.40%  "

Set-TextValue $ws.Range("D20") "471.69"
$ws.Range("E20").Value = "  -3.18%  "

Set-TextValue $ws.Range("D21") "14.43"
$ws.Range("E21").Value = "  -3.68%  "

Set-TextValue $ws.Range("D22") "0.708"
$ws.Range("E22").Value = "  -2.49%  "

Set-TextValue $ws.Range("D23") "7.81"
$ws.Range("E23").Value = "  -0.15%  "

Set-TextValue $ws.Range("D24") "13.47"
$ws.Range("E24").Value = "  -4.15%  "

Set-TextValue $ws.Range("D25") "82.93"
$ws.Range("E25").Value = "  -2.48%  "

$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("E27").Value = "  -6.57%  "

Set-TextValue $ws.Range("D28") "8.40"
$ws.Range("E28").Value = "  -4.75%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D29") "0.119"
$ws.Range("E29").Value = "  -9.84%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D30") "6.96"
$ws.Range("E30").Value = "  +0.21%  "

Set-TextValue $ws.Range("D31") "2.05"
$ws.Range("E31").Value = "  -10.95%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("E33").Value = "  -4.14%  "

Set-TextValue $ws.Range("D34") "26.05"
$ws.Range("E34").Value = "  -4.33%  "

Set-TextValue $ws.Range("D35") "1.11"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").Value = "0.0₃0780"
$ws.Range("E36").Value = "  +4.69%  "

Set-TextValue $ws.Range("D37") "5.91"
$ws.Range("E37").Value = "  -4.52%  "

Set-TextValue $ws.Range("D38") "52.28"
$ws.Range("E38").Value = "  -4.58%  "

Set-TextValue $ws.Range("D39") "453.47"
$ws.Range("E39").Value = "  -5.71%  "

Set-TextValue $ws.Range("D40") "2.96"
$ws.Range("E40").Value = "  -10.71%  "

$ws.Range("E41").Value = "  -4.45%  "

$ws.Range("E42").Value = "  -6.80%  "

Set-TextValue $ws.Range("D43") "8.26"
$ws.Range("E43").Value = "  -3.83%  "

$ws.Range("D44").Value = "2.844.84"
$ws.Range("E44").Value = "  -3.40%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D45") "0.264"
$ws.Range("E45").Value = "  -6.38%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D46") "2.27"
$ws.Range("E46").Value = "  -8.61%  "

$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("E48").Value = "  +0.01%  "

Set-TextValue $ws.Range("D49") "26.04"
$ws.Range("E49").Value = "  -5.97%  "

$ws.Range("E50").Value = "  -3.27%  "

Set-TextValue $ws.Range("D51") "118.62"
$ws.Range("E51").Value = "  -2.05%  "

# Clean up the scratch helper cell so it leaves no trace in the sheet.
$ws.Range("ZZ1").Value = $null
